$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.257.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.691.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -7.18%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -5.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.699.90'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("E12").Value = '  -4.95%  '
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.163.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.334.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.702.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.85%  '
$ws.Range("E18").Value = '  -6.17%  '
$ws.Range("E19").Value = '  -6.71%  '
$ws.Range("E20").Value = '  -6.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '332.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.75%  '
$ws.Range("E22").Value = '  -6.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.995'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.86%  '
$ws.Range("E26").Value = '  -6.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.169'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.996'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0814'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.89%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("E33").Value = '  -5.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.75%  '
$ws.Range("E36").Value = '  -4.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.926'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.64%  '
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.171.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.36%  '
$ws.Range("E42").Value = '  -8.67%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  -4.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.590'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.42%  '
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0224'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.16%  '
$ws.Range("E50").Value = '  -5.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.14%  '
